# Update the 36 placeholder 'Hola' daily-message cells in column C with their
# final unique messages (rows 314-366, skipping the rows that already hold
# image/quote/video triplets).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C314").Value = 'Perquè et mous a dormir (algun cop) a la meva habitació quan hi ha colònies infants!'
$ws.Range("C315").Value = 'Perquè ajudes a tothom de la teva família sempre que pots'
$ws.Range("C317").Value = 'Perquè sempre em preguntes "un iogurt?" quan arribem tard a casa'
$ws.Range("C318").Value = 'Perquè fas stickers de mi SUPER EXTRA xulos'
$ws.Range("C320").Value = 'Perquè m''acompanyes a fer plans MOOOLT guais amb els meus amics (per exemple veure l''Edgar i la Judit participar a soc i sere)'
$ws.Range("C321").Value = 'Perquè quan et faig el dinar dius que esta molt bo encara que sempre siguin raviolis de formatge'
$ws.Range("C323").Value = 'Perquè em recordes els aniversaris de la gent que he de felicitar'
$ws.Range("C324").Value = 'Perquè aguantes les meves bromes de "personatges" que et venen a veure. El meu preferit és el MAGÇAL'
$ws.Range("C325").Value = 'Perquè proves el menjar que fa la meva germana (NORMALMENT MOOOLT DOLENT) i fas veure que esta bo'
$ws.Range("C327").Value = 'Perquè m''acompanyes a fer plans XULISSIMS encara que només m''agradin a mi (BUNQUEER)'
$ws.Range("C328").Value = 'Perquè escoltes els audios super guais que et faig de tant en tant explican-te com em va el dia'
$ws.Range("C330").Value = 'Perquè m''expliques anècdotes molt divertides que et passen a les pràctiques de la universitat'
$ws.Range("C331").Value = 'Perquè em fas regals NINJAS però moolt guais (em vas regalar anar a fer un sushi amb tu a Altafulla pel meu cumple. Quan escric això: 19/05/2025 encara NO ho hem fet. Espero que quan llegeixis això sii jijiji)'
$ws.Range("C334").Value = 'Perquè quan em convides a dinar a casa teva el teu pare fa UN ARROS MOOOOLT BO'
$ws.Range("C335").Value = 'Perquè em dius que em talli el cavell SEMPRE encara que just me''l acabi de tallar'
$ws.Range("C339").Value = 'Perquè m''encanta que (ara que ja portem més de dos anys coneixent-nos) poguem fer plans en dies especials PER SEGON COP'
$ws.Range("C340").Value = 'Perquè SEMPRE estaves preguntant-me si seguiria al taller o no!'
$ws.Range("C341").Value = 'Perquè quan estem els dos junts ens ho passem MOOOLT be'
$ws.Range("C343").Value = 'Perquè sempre intentes buscar un temps per mi encara que estiguem MOOLT ocupats'
$ws.Range("C344").Value = 'Perquè tens una figureta meva de Lego sempre aprop teu a l''habitació'
$ws.Range("C346").Value = 'Perquè em vas regalar un braçelet verd per remplaçar el que se''m va trencar (que l''havies fet tuuu)'
$ws.Range("C347").Value = 'Perquè t''adones de les actualtizacions MEGA GUAIS que faig a la web '
$ws.Range("C348").Value = 'Perquè et vas copiar de la app de la meva agenda MEGA XULA'
$ws.Range("C349").Value = 'Perquè cada dia (o casi cada dia) tens un moment per llegir-te el missatge diari'
$ws.Range("C350").Value = 'Perquè dediques el teu temps en fer regals macos per mi (el mes guai) i els altres'
$ws.Range("C352").Value = 'Perquè em passes algunes fotos de la puça'
$ws.Range("C353").Value = 'Perquè aguantes que li digui MONSTRE a la puça (no és un monstres però es porta BASTANT malament)'
$ws.Range("C355").Value = 'Perquè et preocupes per les meves mans quan estan tallades'
$ws.Range("C356").Value = 'Perquè ets LA MARGARITA més guapa del món (a colònies monis ho vas demostrar)'
$ws.Range("C357").Value = 'Perquè em robes TOTS els mitjons que tinc'
$ws.Range("C359").Value = 'Perquè et preocupes perquè no tingui molta alèrgia quan vinc a casa teva'
$ws.Range("C360").Value = 'Perquè em preguntes sobre com em va la carrera'
$ws.Range("C362").Value = 'Perquè t''has convertit en 2 anys en una persona molt important per mi'
$ws.Range("C364").Value = 'Perquè SEMPRE vols que anem a menjar sushi'
$ws.Range("C365").Value = 'Perquè ets una gran massatgista i vam comprar uns oli de massatges a Altafulla perquè ho poguessis demostrar'
$ws.Range("C366").Value = 'Perquè sempre estàs molt atenta a mi i al que necessito. PREPARA''T QUE DEMÀ ES L''ÚLTIM DIA DE LA WEB. Tot i això, et seguirè estimant sempre.'

# Restore the author's final cursor position / scroll state.
$ws.Range("D140").Select() | Out-Null
